# Corrige os totais de residuos para os anos de 2019 (H) e 2020 (I),
# que estavam divididos por 1000 (provavelmente convertidos de kg para
# toneladas por engano). Os valores sao restaurados para a escala correta
# e a coluna de total_anos (M) e recalculada como soma de B:L em cada linha.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 220
$ws.Range("M3").Value = 705.84
$ws.Range("H5").Value = 255417
$ws.Range("I5").Value = 387522
$ws.Range("M5").Value = 2400728.22
$ws.Range("H7").Value = 3680081
$ws.Range("I7").Value = 3619317
$ws.Range("M7").Value = 29697636.3
$ws.Range("I8").Value = 550
$ws.Range("M8").Value = 991.8499999999999
$ws.Range("I9").Value = 100
$ws.Range("M9").Value = 100
$ws.Range("H10").Value = 37964
$ws.Range("I10").Value = 23563
$ws.Range("M10").Value = 404956.58
$ws.Range("H11").Value = 234539
$ws.Range("I11").Value = 103830
$ws.Range("M11").Value = 4132187.88
$ws.Range("I12").Value = 10
$ws.Range("M12").Value = 97.88
$ws.Range("H13").Value = 200985
$ws.Range("I13").Value = 186372
$ws.Range("M13").Value = 1482007.52
$ws.Range("H14").Value = 65325
$ws.Range("I14").Value = 63206
$ws.Range("M14").Value = 634647.6699999999
$ws.Range("H15").Value = 33518
$ws.Range("I15").Value = 36869
$ws.Range("M15").Value = 271276.08
$ws.Range("I16").Value = 110
$ws.Range("M16").Value = 125.69
$ws.Range("H17").Value = 12092
$ws.Range("I17").Value = 49125
$ws.Range("M17").Value = 61217
$ws.Range("I18").Value = 6750
$ws.Range("M18").Value = 6753
$ws.Range("I20").Value = 89550
$ws.Range("M20").Value = 142750.91
$ws.Range("H21").Value = 49770
$ws.Range("I21").Value = 60075
$ws.Range("M21").Value = 368350.82
$ws.Range("H23").Value = 12984
$ws.Range("I23").Value = 16069
$ws.Range("M23").Value = 81053.85000000001
$ws.Range("I25").Value = 2630
$ws.Range("M25").Value = 2817
$ws.Range("I26").Value = 6810
$ws.Range("M26").Value = 6887
$ws.Range("I27").Value = 74436
$ws.Range("M27").Value = 74436
$ws.Range("H28").Value = 17807
$ws.Range("I28").Value = 40154
$ws.Range("M28").Value = 119319.08
$ws.Range("H29").Value = 218355
$ws.Range("I29").Value = 234256
$ws.Range("M29").Value = 1257367.17
$ws.Range("H30").Value = 447736
$ws.Range("I30").Value = 346544
$ws.Range("M30").Value = 1861233.65
$ws.Range("I31").Value = 5804
$ws.Range("M31").Value = 5804
$ws.Range("H32").Value = 189368
$ws.Range("I32").Value = 174988
$ws.Range("M32").Value = 1507963.54
$ws.Range("H34").Value = 80454
$ws.Range("I34").Value = 79570
$ws.Range("M34").Value = 628428.01
$ws.Range("H35").Value = 81863
$ws.Range("I35").Value = 57967
$ws.Range("M35").Value = 762935.54
$ws.Range("I36").Value = 56530
$ws.Range("M36").Value = 58438.9
$ws.Range("H37").Value = 5629081.1
$ws.Range("I37").Value = 5722927
$ws.Range("M37").Value = 56361609.31999999
